
# Weekly data refresh: a new week's observation is added at the top of the
# data block (row 105) and all existing data rows shift down by one
# (105-119 -> 106-120), so the oldest row (old 119) now lands on row 120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 105:119 down by inserting a new blank row at 105.
$ws.Rows.Item(105).EntireRow.Insert()

# Populate the newly inserted row 105 with this week's record.
$ws.Cells.Item(105, 1).Value = 6
$ws.Cells.Item(105, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(105, 3).Value = "Metropolitana"
$ws.Cells.Item(105, 4).Value = 44504
$ws.Cells.Item(105, 5).Value = 13
$ws.Cells.Item(105, 6).Value = 100112029
$ws.Cells.Item(105, 7).Value = "Orégano"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 35
$ws.Cells.Item(105, 11).Value = 7500
$ws.Cells.Item(105, 12).Value = 8000
$ws.Cells.Item(105, 13).Value = 7729
$ws.Cells.Item(105, 14).Value = "$/docena de atados"
$ws.Cells.Item(105, 15).Value = "Región Metropolitana"
$ws.Cells.Item(105, 16).Value = 2576
$ws.Cells.Item(105, 17).Value = 3
$ws.Cells.Item(105, 18).Value = "Hortaliza"
